$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.805752277374268
$ws.Range("B1").Value = 4.353612899780273
$ws.Range("C1").Value = 3.865099906921387
$ws.Range("D1").Value = 1.429892897605896
$ws.Range("E1").Value = 0.6780727505683899
